$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row texts
$ws.Range("A1").Value = "אימייל"
$ws.Range("D1").Value = "שלב"
$ws.Range("E1").Value = "סטטוס"
$ws.Range("F1").Value = "זמן"

# Update existing row 2 data (was halroy13, becomes candidate0)
$ws.Range("A2").Value = "candidate0@gmail.com"
$ws.Range("B2").Value = "דוד"
$ws.Range("C2").Value = "חי"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "מועמד חדש"
$ws.Range("F2").Value = "2023-01-22 17:09:02.197615+00:00"

# Add new row 3 (candidate1)
$ws.Range("A3").Value = "candidate1@gmail.com"
$ws.Range("B3").Value = "משה"
$ws.Range("C3").Value = "שמחון"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = "מועמד חדש"
$ws.Range("F3").Value = "2023-01-22 17:09:23.363519+00:00"

# Add new row 4 (original halroy13 data, now moved down)
$ws.Range("A4").Value = "halroy13@gmail.com"
$ws.Range("B4").Value = "רוי"
$ws.Range("C4").Value = "הלחמי"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = "מועמד חדש"
$ws.Range("F4").Value = "2023-01-22 16:27:15.458863+00:00"
